$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2222
$ws.Range("I19").Value = 999.8
$ws.Range("J19").Value = 3749.75
$ws.Range("K19").Value = 999.8
$ws.Range("L19").Value = 3749.75
$ws.Range("M19").Value = -824.8
$ws.Range("N19").Value = -4099.75
$ws.Range("H20").Value = 5162
$ws.Range("J20").Value = 10024
$ws.Range("L20").Value = 10024
$ws.Range("N20").Value = -10484
$ws.Range("H35").Value = 5162
$ws.Range("J35").Value = 10024
$ws.Range("L35").Value = 10024
$ws.Range("N35").Value = -10782
$ws.Range("H125").Value = 2000
$ws.Range("I125").Value = 2000
$ws.Range("K125").Value = 18000
$ws.Range("M125").Value = -15540
$ws.Range("H137").Value = 3470.8572
$ws.Range("I137").Value = 3682.8333
$ws.Range("J137").Value = 2199
$ws.Range("K137").Value = 11048.4999
$ws.Range("L137").Value = 6597
$ws.Range("M137").Value = -8498.499899999999
$ws.Range("N137").Value = -11697
$ws.Range("H138").Value = 3829.4736
$ws.Range("J138").Value = 4994.357
$ws.Range("L138").Value = 14983.071
$ws.Range("N138").Value = -25263.071

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1072.8334
$ws.Range("I5").Value = 109.5
$ws.Range("K5").Value = 109.5
$ws.Range("M5").Value = 2.5
$ws.Range("H32").Value = 1194.1786
$ws.Range("I32").Value = 1194.1786
$ws.Range("K32").Value = 1194.1786
$ws.Range("M32").Value = -907.1786
$ws.Range("H45").Value = 4255
$ws.Range("I45").Value = 3299.3333
$ws.Range("J45").Value = 4971.75
$ws.Range("K45").Value = 3299.3333
$ws.Range("L45").Value = 4971.75
$ws.Range("M45").Value = -2922.3333
$ws.Range("N45").Value = -5725.75
$ws.Range("H61").Value = 4612.6665
$ws.Range("J61").Value = 5001.75
$ws.Range("L61").Value = 5001.75
$ws.Range("N61").Value = -5425.75
$ws.Range("N74").Value = ""
$ws.Range("H74").Value = 599.5
$ws.Range("I74").Value = 599.5
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 599.5
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 274.5
$ws.Range("N77").Value = ""
$ws.Range("H77").Value = 599.5
$ws.Range("I77").Value = 599.5
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 2997.5
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = 1370.5
$ws.Range("H102").Value = 962.7778
$ws.Range("I102").Value = 906.75
$ws.Range("K102").Value = 906.75
$ws.Range("M102").Value = 715.25
$ws.Range("H110").Value = 1489.5
$ws.Range("I110").Value = 1600
$ws.Range("K110").Value = 1600
$ws.Range("M110").Value = 445
$ws.Range("H136").Value = 4612.6665
$ws.Range("J136").Value = 5001.75
$ws.Range("L136").Value = 15005.25
$ws.Range("N136").Value = -20105.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1072.8334
$ws.Range("I4").Value = 109.5
$ws.Range("K4").Value = 109.5
$ws.Range("M4").Value = 5.5
$ws.Range("H11").Value = 199.8
$ws.Range("I11").Value = 214.66667
$ws.Range("J11").Value = 177.5
$ws.Range("K11").Value = 214.66667
$ws.Range("L11").Value = 177.5
$ws.Range("M11").Value = -74.66667000000001
$ws.Range("N11").Value = -457.5
$ws.Range("N107").Value = ""
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 873
$ws.Range("I22").Value = 497.33334
$ws.Range("K22").Value = 497.33334
$ws.Range("M22").Value = -147.33334
$ws.Range("N39").Value = ""
$ws.Range("H39").Value = 2000
$ws.Range("I39").Value = 2000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 2000
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -1609
$ws.Range("N49").Value = ""
$ws.Range("H49").Value = 2000
$ws.Range("I49").Value = 2000
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 2000
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -1818
$ws.Range("M58").Value = ""
$ws.Range("N58").Value = ""
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M62").Value = ""
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 2000
$ws.Range("N62").Value = -3248
$ws.Range("M65").Value = ""
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 10000
$ws.Range("N65").Value = -16240
$ws.Range("H134").Value = 13167.286
$ws.Range("I134").Value = 14861.667
$ws.Range("K134").Value = 44585.001
$ws.Range("M134").Value = -42050.001
$ws.Range("M136").Value = ""
$ws.Range("N136").Value = ""
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 83.52941
$ws.Range("J2").Value = 83.125
$ws.Range("L2").Value = 498.75
$ws.Range("N2").Value = -724.75
$ws.Range("H38").Value = 16
$ws.Range("I38").Value = 11
$ws.Range("J38").Value = 18.5
$ws.Range("K38").Value = 33
$ws.Range("L38").Value = 55.5
$ws.Range("M38").Value = 314
$ws.Range("N38").Value = -749.5
$ws.Range("H68").Value = 2406.2222
$ws.Range("J68").Value = 2406.2222
$ws.Range("L68").Value = 7218.6666
$ws.Range("N68").Value = -8840.6666
$ws.Range("H71").Value = 2406.2222
$ws.Range("J71").Value = 2406.2222
$ws.Range("L71").Value = 21655.9998
$ws.Range("N71").Value = -29767.9998

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M31").Value = ""
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M37").Value = ""
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("H80").Value = 810.6667
$ws.Range("I80").Value = 1066.25
$ws.Range("J80").Value = 299.5
$ws.Range("K80").Value = 1066.25
$ws.Range("L80").Value = 299.5
$ws.Range("M80").Value = -68.25
$ws.Range("N80").Value = -2295.5
$ws.Range("H83").Value = 810.6667
$ws.Range("I83").Value = 1066.25
$ws.Range("J83").Value = 299.5
$ws.Range("K83").Value = 5331.25
$ws.Range("L83").Value = 1497.5
$ws.Range("M83").Value = -339.25
$ws.Range("N83").Value = -11481.5
$ws.Range("H132").Value = 4396.8
$ws.Range("J132").Value = 4396.8
$ws.Range("L132").Value = 13190.4
$ws.Range("N132").Value = -18250.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N76").Value = ""
$ws.Range("H76").Value = 1260
$ws.Range("I76").Value = 1260
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 1260
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -922
$ws.Range("N79").Value = ""
$ws.Range("H79").Value = 1260
$ws.Range("I79").Value = 1260
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 1260
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -90
$ws.Range("H132").Value = 3636
$ws.Range("I132").Value = 2999.7144
$ws.Range("K132").Value = 8999.143199999999
$ws.Range("M132").Value = -6469.143199999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 4999.5
$ws.Range("J2").Value = 4999.5
$ws.Range("L2").Value = 4999.5
$ws.Range("N2").Value = -5223.5
$ws.Range("H62").Value = 2995
$ws.Range("I62").Value = 2995
$ws.Range("K62").Value = 2995
$ws.Range("M62").Value = -2371
$ws.Range("H65").Value = 2995
$ws.Range("I65").Value = 2995
$ws.Range("K65").Value = 14975
$ws.Range("M65").Value = -11855
